$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 112083127
$ws.Range("B7").Value = 77724
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 6450
$ws.Range("F7").Value = "Skuggblåslav"
$ws.Range("G7").Value = "Hypogymnia vittata"
$ws.Range("H7").Value = "(Ach.) Parrique"
$ws.Range("Q7").Value = 413052
$ws.Range("R7").Value = 6656343

$ws.Range("A8").Value = 112083125
$ws.Range("B8").Value = 89503
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 5447
$ws.Range("F8").Value = "Vedticka"
$ws.Range("G8").Value = "Fuscoporia viticola"
$ws.Range("H8").Value = "(Schwein.) Murrill"
$ws.Range("Q8").Value = 413016
$ws.Range("R8").Value = 6656415

$ws.Range("A9").Value = 112083118
$ws.Range("B9").Value = 94287
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 53
$ws.Range("F9").Value = "Vedtrappmossa"
$ws.Range("G9").Value = "Crossocalyx hellerianus"
$ws.Range("H9").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q9").Value = 412577
$ws.Range("R9").Value = 6656304

$ws.Range("A10").Value = 112083128
$ws.Range("B10").Value = 77307
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 353
$ws.Range("F10").Value = "Dvärgbägarlav"
$ws.Range("G10").Value = "Cladonia parasitica"
$ws.Range("H10").Value = "(Hoffm.) Hoffm."
$ws.Range("Q10").Value = 413190
$ws.Range("R10").Value = 6656475

$ws.Range("A11").Value = 112083112
$ws.Range("B11").Value = 79566
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 1049
$ws.Range("F11").Value = "Kortskaftad ärgspik"
$ws.Range("G11").Value = "Microcalicium ahlneri"
$ws.Range("H11").Value = "Tibell"
$ws.Range("Q11").Value = 412284
$ws.Range("R11").Value = 6656072

$ws.Range("A12").Value = 112083110
$ws.Range("B12").Value = 78228
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6453
$ws.Range("F12").Value = "Vedskivlav"
$ws.Range("G12").Value = "Hertelidea botryosa"
$ws.Range("H12").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q12").Value = 412206
$ws.Range("R12").Value = 6656051

$ws.Range("A13").Value = 112083126
$ws.Range("B13").Value = 78657
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 229497
$ws.Range("F13").Value = "Korallblylav"
$ws.Range("G13").Value = "Parmeliella triptophylla"
$ws.Range("H13").Value = "(Ach.) Müll.Arg."
$ws.Range("Q13").Value = 413017
$ws.Range("R13").Value = 6656342

$ws.Range("A14").Value = 112083111
$ws.Range("B14").Value = 90800
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = "Dropptaggsvamp"
$ws.Range("G14").Value = "Hydnellum ferrugineum"
$ws.Range("H14").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q14").Value = 412205
$ws.Range("R14").Value = 6655989
